$wb = $excel.ActiveWorkbook
$sprint1 = $wb.Worksheets.Item("Sprint1")
$sprint1.Select()
try { $excel.ActiveWindow.TopLeftCell = $sprint1.Range("E1") ; Write-Output "set TopLeftCell direct" } catch { Write-Output "no direct TopLeftCell prop" }
